$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.884.18"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "2.659.99"
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'521.96"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").Value = "'148.99"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("E7").Value = "  -0.49%  "

$ws.Range("D8").Value = "'0.575"
$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("D9").Value = "2.686.68"
$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").Value = "'6.53"
$ws.Range("E10").Value = "  +2.43%  "

$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").Value = "3.124.55"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").Value = "59.764.29"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").Value = "'21.54"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "2.673.24"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "'348.20"
$ws.Range("E20").Value = "  +0.98%  "

$ws.Range("D21").Value = "'10.67"
$ws.Range("E21").Value = "  +2.15%  "

$ws.Range("E22").Value = "  +2.08%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("D24").Value = "'61.27"
$ws.Range("E24").Value = "  +0.98%  "

$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("D26").Value = "2.774.34"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.163"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'0.992"
$ws.Range("E28").Value = "  -0.87%  "

$ws.Range("D29").Value = "0.0₃0836"
$ws.Range("E29").Value = "  +1.94%  "

$ws.Range("D30").Value = "'7.25"
$ws.Range("E30").Value = "  +2.55%  "

$ws.Range("D31").Value = "'6.63"
$ws.Range("E31").Value = "  +10.53%  "

$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").Value = "'19.15"
$ws.Range("E34").Value = "  +0.82%  "

$ws.Range("D35").Value = "'1.08"
$ws.Range("E35").Value = "  +20.30%  "

$ws.Range("D36").Value = "'149.52"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").Value = "'4.07"
$ws.Range("E37").Value = "  +2.41%  "

$ws.Range("E38").Value = "  +1.71%  "

$ws.Range("D39").Value = "'0.880"
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("D40").Value = "'36.69"
$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("E41").Value = "  +2.30%  "

$ws.Range("D42").Value = "'1.44"
$ws.Range("E42").Value = "  -0.20%  "

$ws.Range("D43").Value = "'290.42"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("D46").Value = "'0.992"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").Value = "'19.80"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("D49").Value = "'4.82"
$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.26"
$ws.Range("E51").Value = "  -1.30%  "
